$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert six new rows before the current row 8 (shifts old rows 8.. down by 6)
# ---------------------------------------------------------------------------
$ws.Rows("8:13").Insert()

# ---------------------------------------------------------------------------
# Copy formatting (number formats / styles) down from the two rows directly
# above the insertion point so the new rows look like the existing data rows.
# Row 6 is a "work session" row (Date/Who/Start/Stop/...), row 7 is the
# "committing to git" summary row that follows it. Row 8/10/12 below mirror
# row 6's layout, rows 9/11/13 mirror row 7's layout.
# ---------------------------------------------------------------------------
$ws.Range("A6:K6").Copy() | Out-Null
$ws.Range("A8:K8").PasteSpecial(-4122) | Out-Null
$ws.Range("A6:K6").Copy() | Out-Null
$ws.Range("A10:K10").PasteSpecial(-4122) | Out-Null
$ws.Range("A6:K6").Copy() | Out-Null
$ws.Range("A12:K12").PasteSpecial(-4122) | Out-Null

$ws.Range("A7:K7").Copy() | Out-Null
$ws.Range("A9:K9").PasteSpecial(-4122) | Out-Null
$ws.Range("A7:K7").Copy() | Out-Null
$ws.Range("A11:K11").PasteSpecial(-4122) | Out-Null
$ws.Range("A7:K7").Copy() | Out-Null
$ws.Range("A13:K13").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Row 8 - new work session (note strings are entered in a specific order so
# that they land in the shared-string table in the desired sequence)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = 41688
$ws.Range("B8").Value = "Doyle"
$ws.Range("C8").Value = 0.75
$ws.Range("D8").Value = 0.79166666666666663
$ws.Range("E8").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("I8").Formula = "=60"
$ws.Range("K8").Value = "Discussed tree of nodes with Steve"

$ws.Range("A9").Value = " =========================    committing to git:"
$ws.Range("B9").Formula = "=A8"
$ws.Range("C9").Formula = "=D8"
$ws.Range("D9").Value = " =========================    "

# ---------------------------------------------------------------------------
# Row 12 next so "Got initial nodes->scene working (no lines etc)" becomes
# shared string index 23, ahead of row 10's string (index 24)
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = 41689
$ws.Range("B12").Value = "Doyle"
$ws.Range("C12").Value = 0.36458333333333331
$ws.Range("D12").Value = 0.39374999999999999
$ws.Range("E12").Value = 0
$ws.Range("G12").Value = 1
$ws.Range("I12").Formula = "=15+27"
$ws.Range("K12").Value = "Got initial nodes->scene working (no lines etc)"

$ws.Range("A13").Value = " =========================    committing to git:"
$ws.Range("B13").Formula = "=A12"
$ws.Range("C13").Formula = "=D12"
$ws.Range("D13").Value = " =========================    "

# ---------------------------------------------------------------------------
# Row 10
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = 41688
$ws.Range("B10").Value = "Doyle"
$ws.Range("C10").Value = 0.92708333333333337
$ws.Range("D10").Value = 0.99375000000000002
$ws.Range("E10").Value = 0
$ws.Range("G10").Value = 1
$ws.Range("I10").Formula = "=45+51"
$ws.Range("K10").Value = "Worked on world to scene"

$ws.Range("A11").Value = " =========================    committing to git:"
$ws.Range("B11").Formula = "=A10"
$ws.Range("C11").Formula = "=D10"
$ws.Range("D11").Value = " =========================    "

# ---------------------------------------------------------------------------
# Update the selection to match the edited workbook
# ---------------------------------------------------------------------------
$ws.Range("A4:XFD13").Select()

Write-Host "done"
